$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HH")

# New header cells for VAT and Gross columns.
$ws.Range("J11").Value = "VAT"
$ws.Range("K11").Value = "Gross"

# New VAT / Gross figures for the two activity rows.
$ws.Range("J12").Value = 13.53
$ws.Range("K12").Value = 150.56
$ws.Range("J13").Value = 13.53
$ws.Range("K13").Value = 150.56

# C13's "Cop" cell picks up the plain/general number format shared by the
# rest of the column (matches the rest of the parser's default styling).
$ws.Range("C13").NumberFormat = $ws.Range("D12").NumberFormat()

# Leave the selection where the author's last edit was.
$null = $ws.Range("K12").Select()
